$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 295, shifting the existing
# rows 295-305 down to 297-307 (formatting of row above is copied along,
# which preserves the date style "s=2" on column D).
$ws.Range("A295:T296").EntireRow.Insert()

# Populate the two newly inserted rows (295 and 296) with this week's data.
# Columns that are constant across the whole data block (A, B, C, E, F, G,
# H, I, J, K, T) are copied from the surrounding rows.

# Row 295
$ws.Range("A295").Value = 10
$ws.Range("B295").Value = "Vega Modelo de Temuco"
$ws.Range("C295").Value = "La Araucanía"
$ws.Range("D295").Value = 44931
$ws.Range("E295").Value = 9
$ws.Range("F295").Value = "Fruta"
$ws.Range("G295").Value = 100101
$ws.Range("H295").Value = "Berries"
$ws.Range("I295").Value = 100112025
$ws.Range("J295").Value = "Frutilla"
$ws.Range("K295").Value = "Sin especificar"
$ws.Range("L295").Value = "Primera"
$ws.Range("M295").Value = 40
$ws.Range("N295").Value = 7000
$ws.Range("O295").Value = 7000
$ws.Range("P295").Value = 7000
$ws.Range("Q295").Value = "$/caja 7 kilos"
$ws.Range("R295").Value = "Región de La Araucanía"
$ws.Range("S295").Value = 1000
$ws.Range("T295").Value = 7

# Row 296
$ws.Range("A296").Value = 10
$ws.Range("B296").Value = "Vega Modelo de Temuco"
$ws.Range("C296").Value = "La Araucanía"
$ws.Range("D296").Value = 44931
$ws.Range("E296").Value = 9
$ws.Range("F296").Value = "Fruta"
$ws.Range("G296").Value = 100101
$ws.Range("H296").Value = "Berries"
$ws.Range("I296").Value = 100112025
$ws.Range("J296").Value = "Frutilla"
$ws.Range("K296").Value = "Sin especificar"
$ws.Range("L296").Value = "Segunda"
$ws.Range("M296").Value = 50
$ws.Range("N296").Value = 4000
$ws.Range("O296").Value = 5000
$ws.Range("P296").Value = 4600
$ws.Range("Q296").Value = "$/caja 7 kilos"
$ws.Range("R296").Value = "Región de La Araucanía"
$ws.Range("S296").Value = 657
$ws.Range("T296").Value = 7
